$d = $word.ActiveDocument

# --- Body (document.xml): "QWREW" -> "QWR" (bold run referring to the school/DE name) ---
$body = $d.Content
$body.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)

# --- Header (header1.xml): sequence of placeholder replacements ---
$hdr = $d.Sections.Item(1).Headers.Item(1)

# "DIRETORIA DE ENSINO REGIAO REW" -> "... QWER"
$rng = $hdr.Range.Duplicate()
$rng.Find.Execute("REW", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 2)

# "QWREW - DEP." -> "QWR - DEP."
$rng = $hdr.Range.Duplicate()
$rng.Find.Execute("QWREW", $true, $true, $false, $false, $false, $true, 1, $false, "QWR", 2)

# Address line runs "Rew" (x5) -> "Qwer"
$rng = $hdr.Range.Duplicate()
$rng.Find.Execute("Rew", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 2)

# CEP / Tel / Email values "rew" (x3) -> "qwer"
$rng = $hdr.Range.Duplicate()
$rng.Find.Execute("rew", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 2)
